# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values recalculated for each row (A1:J26 sheet, K is column G)
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
